$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.993.73'
$ws.Range('E2').Value = '  +4.87%  '

$ws.Range('D3').Value = '2.778.36'
$ws.Range('E3').Value = '  +4.97%  '

$ws.Range('E4').Value = '  +0.03%  '

$ws.Range('B5').Value = 'Solana'
$ws.Range('C5').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D5').Value = "'115.30"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.37%  '

$ws.Range('B6').Value = 'BNB'
$ws.Range('C6').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D6').Value = "'338.94"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.76%  '

$ws.Range('D7').Value = "'0.545"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.48%  '

$ws.Range('E8').Value = '  -0.03%  '

$ws.Range('D9').Value = "'0.575"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.67%  '

$ws.Range('D10').Value = "'41.70"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.04%  '

$ws.Range('D11').Value = "'0.0858"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +5.29%  '

$ws.Range('D12').Value = "'20.08"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.24%  '

$ws.Range('E13').Value = '  +1.77%  '

$ws.Range('D14').Value = "'7.58"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.16%  '

$ws.Range('D15').Value = '3.210.17'
$ws.Range('E15').Value = '  +4.92%  '

$ws.Range('D16').Value = '2.781.51'
$ws.Range('E16').Value = '  +5.31%  '

$ws.Range('D17').Value = '51.816.38'
$ws.Range('E17').Value = '  +4.54%  '

$ws.Range('D18').Value = "'0.876"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.88%  '

$ws.Range('D19').Value = "'3.19"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +9.97%  '

$ws.Range('D20').Value = "'6.97"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +4.38%  '

$ws.Range('D21').Value = "'13.19"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.41%  '

$ws.Range('D22').Value = '0.0₃0976'
$ws.Range('E22').Value = '  +2.90%  '

$ws.Range('D23').Value = "'275.90"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.79%  '

$ws.Range('D24').Value = "'69.85"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.45%  '

$ws.Range('D25').Value = "'2.73"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +6.07%  '

$ws.Range('D26').Value = "'26.67"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.17%  '

$ws.Range('D28').Value = "'10.15"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.10%  '

$ws.Range('E29').Value = '  +0.88%  '

$ws.Range('E30').Value = '  +1.58%  '

$ws.Range('D31').Value = "'34.60"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.20%  '

$ws.Range('D32').Value = "'50.16"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.08%  '

$ws.Range('D33').Value = "'5.69"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.08%  '

$ws.Range('D34').Value = "'0.0818"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.32%  '

$ws.Range('E35').Value = '  -0.14%  '

$ws.Range('D36').Value = "'18.95"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.42%  '

$ws.Range('E37').Value = '  +2.88%  '

$ws.Range('D38').Value = "'4.93"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.17%  '

$ws.Range('D39').Value = "'3.20"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.54%  '

$ws.Range('E40').Value = '  +7.79%  '

$ws.Range('D41').Value = "'2.67"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +26.19%  '

$ws.Range('E42').Value = '  +3.58%  '

$ws.Range('E43').Value = '  +3.00%  '

$ws.Range('D44').Value = "'125.65"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.27%  '

$ws.Range('D45').Value = "'23.10"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.33%  '

$ws.Range('D46').Value = '2.066.66'
$ws.Range('E46').Value = '  +0.14%  '

$ws.Range('D47').Value = "'3.30"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.07%  '

$ws.Range('E48').Value = '  +1.16%  '

$ws.Range('D49').Value = "'5.54"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +5.65%  '

$ws.Range('D50').Value = "'8.87"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.70%  '

$ws.Range('B51').Value = 'SEI'
$ws.Range('C51').Value = 'https://coinranking.com/coin/8nxCqs-uj+sei-sei'
$ws.Range('D51').Value = "'0.877"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +15.26%  '
